$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A13").Value = "Minimum Number Of Dimensions For Output Tensor"
$ws.Range("B16").Value = "Immediate Bug Fix. If Paid User, Then Feature And Bug Fix Requests From Them Be Prioritized First."
$ws.Range("C16").Value = "Immediate Bug Fix. If Paid User, Then Feature And Bug Fix Requests From Them Be Prioritized First."

$ws.Range("F24").Select()
